$d = $word.ActiveDocument

# The resume header paragraph currently contains only the name
# ("Dheeraj Chand"). The short-resume template is missing the contact
# info line that the long-resume template has, so insert a new,
# center-aligned paragraph directly after the name with the phone,
# email, website, LinkedIn and location, matching the long-resume
# formatting (no bold / no explicit run formatting).
$d.Content.Find.Execute(
    "Dheeraj Chand",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
